$wb = $excel.ActiveWorkbook

$wsSource = $wb.Worksheets.Item("source")
$wsZh     = $wb.Worksheets.Item("ZH test")

# --- 1. Populate 'ZH test' rows 20:38 with the same strings as rows 1:19 -----
# (duplicate the 19 existing Chinese lines so the 'source' sheet's formulas
#  have real text to resolve to instead of falling through to 0)
for ($i = 1; $i -le 19; $i++) {
    $srcCell = $wsZh.Cells.Item($i, 1)
    $dstCell = $wsZh.Cells.Item($i + 19, 1)
    $dstCell.Value = $srcCell.Text
}

# Reflect the edit as the user would have seen it: scroll/select the newly
# typed block on the 'ZH test' sheet.
$wsZh.Activate()
$wsZh.Range("A20:A38").Select()

# --- 2. Switch the language selector on 'source' to Chinese ------------------
$wsSource.Activate()
$wsSource.Range("B1").Value = "中文"

# --- 3. Drop the trailing placeholder rows (35:46) that only ever evaluated
#        to 0 (out of range of every language sheet) ---------------------------
$wsSource.Range("A35:A46").EntireRow.Delete()

# Leave the selection/scroll position the way the author left it.
$wsSource.Range("A35:A47").Select()
